{"js": "// fix: update by user require 15/03/2024\n//\n// 1) In the \"Th\u1eddi h\u1ea1n\" (loan term) row, the run that shows the literal\n//    word \"th\u00e1ng\" (month) is still in the legacy \".VnTime\" font. Re-font\n//    it to \"Times New Roman\" and extend the sentence so it reads\n//    \"{month_count} th\u00e1ng ({period_count} k\u1ef3)\" \u2014 i.e. add a leading\n//    space (styled like the maroon {month_count} placeholder right\n//    before it) and a trailing \" ({period_count} k\u1ef3)\" suffix.\n// 2) At the very end of the document (signature block) append the\n//    \"{user_full_name}\" placeholder right after the trailing spaces run.\n\n// ---- 1) \"th\u00e1ng\" term run -------------------------------------------------\nconst searchResults = context.document.body.search(\"th\u00e1ng\", { matchCase: true });\nsearchResults.load(\"items/font/name\");\nawait context.sync();\n\nlet thangRange = null;\nfor (const r of searchResults.items) {\n  if (r.font.name === \".VnTime\") {\n    thangRange = r;\n    break;\n  }\n}\n\nif (thangRange) {\n  // Re-font the existing \"th\u00e1ng\" run to Times New Roman.\n  thangRange.font.set({ name: \"Times New Roman\" });\n\n  // Append \" ({period_count} k\u1ef3)\" right after \"th\u00e1ng\" (same style).\n  const suffixRange = thangRange.insertText(\" ({period_count} k\u1ef3)\", Word.InsertLocation.after);\n  suffixRange.font.set({ name: \"Times New Roman\" });\n\n  // Insert a maroon space before \"th\u00e1ng\" (matches the {month_count} run).\n  const prefixRange = thangRange.insertText(\" \", Word.InsertLocation.before);\n  prefixRange.font.set({ name: \"Times New Roman\", color: \"#800000\" });\n\n  await context.sync();\n}\n\n// ---- 2) append \"{user_full_name}\" at the end of the document -------------\nconst paragraphs = context.document.body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nconst lastParagraph = paragraphs.items[paragraphs.items.length - 1];\nconst nameRange = lastParagraph.insertText(\"{user_full_name}\", Word.InsertLocation.end);\nnameRange.font.set({ name: \"Times New Roman\" });\n\nawait context.sync();\n", "ps1": "# fix: update by user require 15/03/2024\n#\n# 1) In the \"Th\u1eddi h\u1ea1n\" (loan term) row, the run that shows the literal\n#    word \"th\u00e1ng\" (month) is still in the legacy \".VnTime\" font. Re-font\n#    it to \"Times New Roman\" and extend the sentence so it reads\n#    \"{month_count} th\u00e1ng ({period_count} k\u1ef3)\" -- i.e. add a leading\n#    space (styled like the maroon {month_count} placeholder right\n#    before it) and a trailing \" ({period_count} k\u1ef3)\" suffix.\n# 2) At the very end of the document (signature block) append the\n#    \"{user_full_name}\" placeholder right after the trailing spaces run.\n\n$d = $word.ActiveDocument\n\n# ---- 1) \"th\u00e1ng\" term run ---------------------------------------------------\n$scan = $d.Content\n$found = $scan.Find\n$found.ClearFormatting()\n$found.Text = \"th\u00e1ng\"\n$found.MatchCase = $true\n$found.MatchWholeWord = $false\n$found.Forward = $true\n$found.Wrap = 0\n\n$target = $null\nwhile ($found.Execute()) {\n  if ($scan.Font.Name -eq \".VnTime\") {\n    $target = $d.Range($scan.Start, $scan.End)\n    break\n  }\n  $scan.Collapse(0)\n}\n\nif ($target -ne $null) {\n  $origStart = $target.Start\n  $origEnd = $target.End\n\n  # Re-font the existing \"th\u00e1ng\" run to Times New Roman.\n  $target.Font.Name = \"Times New Roman\"\n\n  # Append \" ({period_count} k\u1ef3)\" right after \"th\u00e1ng\" (same style).\n  $suffixRng = $d.Range($origEnd, $origEnd)\n  $suffixRng.InsertAfter(\" ({period_count} k\u1ef3)\")\n  $suffixRng.Font.Name = \"Times New Roman\"\n\n  # Insert a maroon space before \"th\u00e1ng\" (matches the {month_count} run).\n  $prefixRng = $d.Range($origStart, $origStart)\n  $prefixRng.InsertBefore(\" \")\n  $prefixRng.Font.Name = \"Times New Roman\"\n  $prefixRng.Font.Color = 128\n}\n\n# ---- 2) append \"{user_full_name}\" at the end of the document --------------\n$d.Content.InsertAfter(\"{user_full_name}\")\n\n$endScan = $d.Content\nif ($endScan.Find.Execute(\"{user_full_name}\")) {\n  $endScan.Font.Name = \"Times New Roman\"\n}\n"}
